$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Text -eq $oldVal) {
        $cell.Value = $newVal
    }
}
